$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
function Set-TextCell([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "29.164.91"
Set-TextCell "E2" "  +0.31%  "
Set-TextCell "D3" "1.826.13"
Set-TextCell "E3" "  -0.10%  "
Set-TextCell "D4" "0.9996"
Set-TextCell "E4" "  +0.16%  "
Set-TextCell "D5" "241.68"
Set-TextCell "E5" "  -0.49%  "
Set-TextCell "D6" "0.6188"
Set-TextCell "E6" "  -0.53%  "
Set-TextCell "D7" "1.001"
Set-TextCell "E7" "  -0.02%  "
Set-TextCell "D8" "0.07325"
Set-TextCell "E8" "  -1.46%  "
Set-TextCell "D9" "0.2888"
Set-TextCell "E9" "  -1.04%  "
Set-TextCell "D10" "23.02"
Set-TextCell "E10" "  -0.74%  "
Set-TextCell "D11" "0.07679"
Set-TextCell "E11" "  -0.22%  "
Set-TextCell "D12" "1.826.02"
Set-TextCell "E12" "  -0.30%  "
Set-TextCell "D13" "4.957"
Set-TextCell "E13" "  -0.91%  "
Set-TextCell "D14" "0.6637"
Set-TextCell "E14" "  -0.47%  "
Set-TextCell "D15" "82.25"
Set-TextCell "E15" "  -0.21%  "
Set-TextCell "D16" "0.000008934"
Set-TextCell "E16" "  -4.78%  "
Set-TextCell "D17" "5.858"
Set-TextCell "E17" "  -1.52%  "
Set-TextCell "D18" "29.133.04"
Set-TextCell "E18" "  +0.23%  "
Set-TextCell "D19" "2.070.47"
Set-TextCell "E19" "  +0.01%  "
Set-TextCell "D20" "237.93"
Set-TextCell "E20" "  +6.92%  "
Set-TextCell "D21" "12.42"
Set-TextCell "E21" "  -1.22%  "
Set-TextCell "E22" "  -0.06%  "
Set-TextCell "D23" "7.284"
Set-TextCell "E23" "  +2.54%  "
Set-TextCell "E24" "  +0.09%  "
Set-TextCell "D25" "158.24"
Set-TextCell "E25" "  -1.01%  "
Set-TextCell "D26" "0.1423"
Set-TextCell "E26" "  +2.46%  "
Set-TextCell "D27" "8.479"
Set-TextCell "E27" "  -0.04%  "
Set-TextCell "D28" "17.66"
Set-TextCell "E28" "  -1.00%  "
Set-TextCell "D29" "1.480"
Set-TextCell "E29" "  -0.55%  "
Set-TextCell "D30" "0.05599"
Set-TextCell "E30" "  -2.88%  "
Set-TextCell "D31" "4.084"
Set-TextCell "E31" "  -0.93%  "
Set-TextCell "D32" "4.093"
Set-TextCell "E32" "  -1.45%  "
Set-TextCell "E33" "  -0.58%  "
Set-TextCell "D34" "1.837"
Set-TextCell "E34" "  +0.39%  "
Set-TextCell "D35" "0.7333"
Set-TextCell "E35" "  -0.55%  "
Set-TextCell "D36" "1.130"
Set-TextCell "E36" "  -0.54%  "
Set-TextCell "D37" "2.626"
Set-TextCell "E37" "  -1.58%  "
Set-TextCell "D38" "2.843"
Set-TextCell "E38" "  +3.02%  "
Set-TextCell "D39" "1.212.93"
Set-TextCell "E39" "  -0.68%  "
Set-TextCell "D40" "0.01766"
Set-TextCell "E40" "  -0.11%  "
Set-TextCell "D41" "6.304"
Set-TextCell "E41" "  -2.67%  "
Set-TextCell "D42" "0.9204"
Set-TextCell "E42" "  +3.44%  "
Set-TextCell "D43" "1.001"
Set-TextCell "E43" "  +0.04%  "
Set-TextCell "D44" "101.69"
Set-TextCell "E44" "  -0.41%  "
Set-TextCell "D45" "1.972.48"
Set-TextCell "E45" "  -0.35%  "
Set-TextCell "D46" "64.82"
Set-TextCell "E46" "  -1.58%  "
Set-TextCell "D47" "0.5090"
Set-TextCell "E47" "  +0.18%  "
Set-TextCell "B48" "TheSandbox"
Set-TextCell "C48" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D48" "0.4024"
Set-TextCell "E48" "  -0.56%  "
Set-TextCell "B49" "EnergySwap"
Set-TextCell "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D49" "9.106"
Set-TextCell "E49" "  +1.42%  "
Set-TextCell "B50" "BabyDogeCoin"
Set-TextCell "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D50" "0.00000000116"
Set-TextCell "E50" "  -6.60%  "
Set-TextCell "D51" "0.05760"
Set-TextCell "E51" "  -1.09%  "
